$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.142.67"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.001.14"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.52"
$ws.Range("E5").Value = "  +8.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.58"
$ws.Range("E6").Value = "  +7.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.686"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.751"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.38"
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.98"
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.635.53"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.004.30"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("E16").Value = "  +8.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.42"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.861.58"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.45"
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.80"
$ws.Range("E22").Value = "  +13.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.56"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.29"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.44"
$ws.Range("E26").Value = "  +17.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.96"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.47"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.80"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.131"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "49.27"
$ws.Range("E34").Value = "  -4.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "673.02"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.90"
$ws.Range("E36").Value = "  +9.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.443"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("E38").Value = "  +6.71%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0490"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.79"
$ws.Range("E45").Value = "  +9.74%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.838.70"
$ws.Range("E50").Value = "  +12.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.39"
$ws.Range("E51").Value = "  +4.33%  "
